$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.096.02"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -4.46%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.097.93"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -4.88%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "603.10"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.73%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.24"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -8.89%  "
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.096.00"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -4.92%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.516"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -4.50%  "
$ws.Range("E10").Value = "  -7.42%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.19"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -9.13%  "
$ws.Range("E12").Value = "  -5.96%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000248"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -7.85%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.93"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -9.33%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.608.17"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -4.81%  "
$ws.Range("E16").Value = "  +1.45%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.211.19"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -4.39%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.103.12"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -4.87%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.74"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -7.71%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "470.36"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -5.56%  "
$ws.Range("E21").Value = "  -5.47%  "
$ws.Range("E22").Value = "  -6.73%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.65"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -5.01%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.36"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -8.44%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "82.98"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -4.20%  "
$ws.Range("E26").Value = "  +0.09%  "
$ws.Range("E27").Value = "  -8.90%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.31"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -8.88%  "
$ws.Range("B29").Value = "NEARProtocol"
$ws.Range("C29").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.80"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -3.36%  "
$ws.Range("B30").Value = "Hedera"
$ws.Range("C30").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.116"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -11.35%  "
$ws.Range("B31").Value = "ImmutableX"
$ws.Range("C31").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.06"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -12.50%  "
$ws.Range("E32").Value = "  -0.01%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.63"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -7.25%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "26.02"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -6.40%  "
$ws.Range("E35").Value = "  -3.39%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.87"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -8.08%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "52.38"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -5.74%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0₃0740"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -5.84%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "452.02"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -8.75%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.90"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -16.22%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0388"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -7.98%  "
$ws.Range("E42").Value = "  -9.95%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.26"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -5.68%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.826.40"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -5.67%  "
$ws.Range("E45").Value = "  -9.83%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.24"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -12.60%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.40"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -2.56%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "25.80"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -10.05%  "
$ws.Range("E50").Value = "  -5.65%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "118.16"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -2.51%  "
